$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Bug fix: spaces are possible when registering -> mark "Ratings can be
# assigned to Exemplar by user" requirement as implemented ("x" instead of "tbd")
$ws.Range("B11").Value = "x"

# Update the active selection left over from editing, as reflected in the saved file
$ws.Range("G14").Select()
